$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column G ("K") values for rows 2-5 (Strike# -> K regen)
$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 6
